$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B2 value from 2 to 1
$ws.Range("B2").Value = 1

# Copy formatting from A2 (which has style s="1") to A3
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats

# Set new row 3 values: A3 = 2, B3 = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
